# Updated cryptos list - applies Price (D) and Volume(1h) (E) updates
# for each coin row, matching the latest scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.921.57'
$ws.Range("E2").Value = '  -0.79%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.862.05'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.78'
$ws.Range("E5").Value = '  -0.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5060'
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3629'
$ws.Range("E8").Value = '  -3.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07161'
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8941'
$ws.Range("E10").Value = '  +0.39%  '
$ws.Range("E11").Value = '  -0.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.854.54'
$ws.Range("E12").Value = '  -0.89%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07443'
$ws.Range("E13").Value = '  -1.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.86'
$ws.Range("E14").Value = '  +3.85%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.232'
$ws.Range("E15").Value = '  -1.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  -0.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008475'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.10'
$ws.Range("E18").Value = '  -0.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.966.65'
$ws.Range("E20").Value = '  -0.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.018'
$ws.Range("E21").Value = '  -1.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.089.53'
$ws.Range("E22").Value = '  -1.65%  '
$ws.Range("E23").Value = '  -2.86%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.433'
$ws.Range("E24").Value = '  -0.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.13'
$ws.Range("E25").Value = '  -1.99%  '
$ws.Range("E26").Value = '  -2.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.86'
$ws.Range("E27").Value = '  -0.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.065'
$ws.Range("E28").Value = '  -1.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.24'
$ws.Range("E29").Value = '  +0.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.669'
$ws.Range("E30").Value = '  -2.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.675'
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09230'
$ws.Range("E32").Value = '  +2.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05084'
$ws.Range("E33").Value = '  -1.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.997'
$ws.Range("E34").Value = '  -3.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7457'
$ws.Range("E36").Value = '  -1.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.288'
$ws.Range("E37").Value = '  +8.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.522'
$ws.Range("E38").Value = '  -1.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01999'
$ws.Range("E39").Value = '  -1.75%  '
$ws.Range("E40").Value = '  +0.78%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5364'
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '117.96'
$ws.Range("E42").Value = '  +2.39%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.494'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.528'
$ws.Range("E44").Value = '  +0.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1468'
$ws.Range("E45").Value = '  -0.82%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4653'
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9990'
$ws.Range("E47").Value = '  -0.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.05'
$ws.Range("E48").Value = '  -0.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.562'
$ws.Range("E49").Value = '  -0.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.88'
$ws.Range("E50").Value = '  +0.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.10'
$ws.Range("E51").Value = '  -2.70%  '
